# Updates the Price (D) and Volume(1h) (E) columns on the crypto list
# to reflect the latest values scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = <new price text or $null if unchanged>; E = <new volume text> }
$updates = @{
    2 = @{ D = '25.781.68'; E = '  +0.43%  ' }
    3 = @{ D = '1.747.92'; E = '  +0.15%  ' }
    4 = @{ D = '1.002'; E = '  +0.03%  ' }
    5 = @{ D = '235.76'; E = '  -0.40%  ' }
    6 = @{ D = $null; E = '  -0.01%  ' }
    7 = @{ D = '0.5079'; E = '  +3.56%  ' }
    8 = @{ D = '40.58'; E = '  -2.36%  ' }
    9 = @{ D = '0.2672'; E = '  +7.57%  ' }
    10 = @{ D = '0.06181'; E = '  +3.57%  ' }
    11 = @{ D = '1.751.68'; E = '  +0.46%  ' }
    12 = @{ D = '0.06930'; E = '  +2.05%  ' }
    13 = @{ D = '15.42'; E = '  +4.63%  ' }
    14 = @{ D = '0.6252'; E = '  +11.54%  ' }
    15 = @{ D = '4.468'; E = '  +0.27%  ' }
    16 = @{ D = '77.55'; E = '  +0.60%  ' }
    17 = @{ D = $null; E = '  +0.01%  ' }
    18 = @{ D = $null; E = '  -0.03%  ' }
    19 = @{ D = '25.798.65'; E = '  +0.34%  ' }
    20 = @{ D = '11.62'; E = '  +1.89%  ' }
    21 = @{ D = '0.000006658'; E = '  +1.92%  ' }
    22 = @{ D = '1.976.90'; E = '  +0.56%  ' }
    23 = @{ D = '4.053'; E = '  +2.10%  ' }
    24 = @{ D = '8.253'; E = '  +5.24%  ' }
    25 = @{ D = '5.128'; E = '  +2.47%  ' }
    26 = @{ D = '136.72'; E = '  +0.30%  ' }
    27 = @{ D = $null; E = '  -2.06%  ' }
    28 = @{ D = '15.09'; E = '  +3.25%  ' }
    29 = @{ D = '1.742'; E = '  -3.44%  ' }
    30 = @{ D = '102.46'; E = '  +0.80%  ' }
    31 = @{ D = '0.08184'; E = '  +2.22%  ' }
    32 = @{ D = '3.695'; E = '  -1.47%  ' }
    33 = @{ D = '3.396'; E = '  +2.94%  ' }
    34 = @{ D = '0.04414'; E = '  +0.72%  ' }
    35 = @{ D = '2.656'; E = '  +2.81%  ' }
    36 = @{ D = '0.9950'; E = '  +1.55%  ' }
    37 = @{ D = '0.6000'; E = '  -0.79%  ' }
    38 = @{ D = '2.623'; E = '  -1.93%  ' }
    39 = @{ D = '0.01563'; E = '  +4.63%  ' }
    40 = @{ D = '1.937'; E = '  -3.08%  ' }
    41 = @{ D = $null; E = '  +0.02%  ' }
    42 = @{ D = '101.32'; E = '  -1.95%  ' }
    43 = @{ D = '0.7506'; E = '  -0.66%  ' }
    44 = @{ D = '0.3822'; E = '  +3.55%  ' }
    45 = @{ D = '4.892'; E = '  -4.80%  ' }
    46 = @{ D = '0.05500'; E = '  +6.76%  ' }
    47 = @{ D = '0.1096'; E = '  +2.55%  ' }
    48 = @{ D = '5.936'; E = '  +1.49%  ' }
    49 = @{ D = '30.02'; E = '  +0.09%  ' }
    50 = @{ D = '52.57'; E = '  +0.44%  ' }
    51 = @{ D = $null; E = '  +0.57%  ' }
}

foreach ($row in $updates.Keys) {
    $u = $updates[$row]
    if ($null -ne $u.D) {
        $cell = $ws.Range("D$row")
        # Force text so numeric-looking prices (e.g. "1.002") are not
        # reinterpreted as numbers; restore the default style afterwards
        # so the cell keeps its original (unstyled) appearance.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    $ws.Range("E$row").Value = $u.E
}
